$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - "time_taken", using the same bold/centered header style
# as the other header cells (B1:E1). Copy E1's formatting over first (this
# brings the value across too), then overwrite the value.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Timestamps for F2:F59, one per data row, matching the commit's metadata values.
$timestamps = @(
    "2021-10-05 13:39:13.867719",
    "2021-10-05 13:39:13.867731",
    "2021-10-05 13:39:13.867735",
    "2021-10-05 13:39:13.867738",
    "2021-10-05 13:39:13.867740",
    "2021-10-05 13:39:13.867743",
    "2021-10-05 13:39:13.867746",
    "2021-10-05 13:39:13.867748",
    "2021-10-05 13:39:13.867751",
    "2021-10-05 13:39:13.867754",
    "2021-10-05 13:39:13.867756",
    "2021-10-05 13:39:13.867759",
    "2021-10-05 13:39:13.867762",
    "2021-10-05 13:39:13.867764",
    "2021-10-05 13:39:13.867767",
    "2021-10-05 13:39:13.867770",
    "2021-10-05 13:39:13.867772",
    "2021-10-05 13:39:13.867775",
    "2021-10-05 13:39:13.867778",
    "2021-10-05 13:39:13.867781",
    "2021-10-05 13:39:13.867783",
    "2021-10-05 13:39:13.867786",
    "2021-10-05 13:39:13.867789",
    "2021-10-05 13:39:13.867791",
    "2021-10-05 13:39:13.867794",
    "2021-10-05 13:39:13.867797",
    "2021-10-05 13:39:13.867800",
    "2021-10-05 13:39:13.867803",
    "2021-10-05 13:39:13.867805",
    "2021-10-05 13:39:13.867808",
    "2021-10-05 13:39:13.867810",
    "2021-10-05 13:39:13.867813",
    "2021-10-05 13:39:13.867816",
    "2021-10-05 13:39:13.867818",
    "2021-10-05 13:39:13.867821",
    "2021-10-05 13:39:13.867823",
    "2021-10-05 13:39:13.867826",
    "2021-10-05 13:39:13.867829",
    "2021-10-05 13:39:13.867831",
    "2021-10-05 13:39:13.867834",
    "2021-10-05 13:39:13.867837",
    "2021-10-05 13:39:13.867839",
    "2021-10-05 13:39:13.867842",
    "2021-10-05 13:39:13.867844",
    "2021-10-05 13:39:13.867847",
    "2021-10-05 13:39:13.867849",
    "2021-10-05 13:39:13.867852",
    "2021-10-05 13:39:13.867854",
    "2021-10-05 13:39:13.867857",
    "2021-10-05 13:39:13.867860",
    "2021-10-05 13:39:13.867862",
    "2021-10-05 13:39:13.867865",
    "2021-10-05 13:39:13.867868",
    "2021-10-05 13:39:13.867870",
    "2021-10-05 13:39:13.867873",
    "2021-10-05 13:39:13.867875",
    "2021-10-05 13:39:13.867878",
    "2021-10-05 13:39:13.867880"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
